# Auto-generated edit script applying numeric updates per the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 6429222.5
$ws.Range("J17").Value = 7500660
$ws.Range("L17").Value = 22501980
$ws.Range("N17").Value = -22502316
$ws.Range("H18").Value = 429.5
$ws.Range("I18").Value = 429.5
$ws.Range("K18").Value = 429.5
$ws.Range("M18").Value = -145.5
$ws.Range("H58").Value = 958.1818
$ws.Range("J58").Value = 2250
$ws.Range("L58").Value = 6750
$ws.Range("N58").Value = -7050
$ws.Range("H64").Value = 6633.6665
$ws.Range("I64").Value = 6640.6
$ws.Range("J64").Value = 6599
$ws.Range("K64").Value = 6640.6
$ws.Range("L64").Value = 6599
$ws.Range("M64").Value = -6392.6
$ws.Range("N64").Value = -7095
$ws.Range("H67").Value = 6633.6665
$ws.Range("I67").Value = 6640.6
$ws.Range("J67").Value = 6599
$ws.Range("K67").Value = 6640.6
$ws.Range("L67").Value = 6599
$ws.Range("M67").Value = -5782.6
$ws.Range("N67").Value = -8315
$ws.Range("H76").Value = 6642.7144
$ws.Range("I76").Value = 5428.2856
$ws.Range("J76").Value = 7857.143
$ws.Range("K76").Value = 5428.2856
$ws.Range("L76").Value = 7857.143
$ws.Range("M76").Value = -5113.2856
$ws.Range("N76").Value = -8487.143
$ws.Range("H79").Value = 6642.7144
$ws.Range("I79").Value = 5428.2856
$ws.Range("J79").Value = 7857.143
$ws.Range("K79").Value = 5428.2856
$ws.Range("L79").Value = 7857.143
$ws.Range("M79").Value = -4336.2856
$ws.Range("N79").Value = -10041.143
$ws.Range("H86").Value = 15155251
$ws.Range("I86").Value = 2928.7
$ws.Range("K86").Value = 2928.7
$ws.Range("M86").Value = -1805.7
$ws.Range("H89").Value = 15155251
$ws.Range("I89").Value = 2928.7
$ws.Range("K89").Value = 14643.5
$ws.Range("M89").Value = -9027.5
$ws.Range("H113").Value = 6950.5
$ws.Range("I113").Value = 6833.3335
$ws.Range("K113").Value = 6833.3335
$ws.Range("M113").Value = -3579.3335
$ws.Range("H114").Value = 27600
$ws.Range("J114").Value = 27600
$ws.Range("L114").Value = 27600
$ws.Range("N114").Value = -36278
$ws.Range("H129").Value = 2023.6428
$ws.Range("I129").Value = 1820.8
$ws.Range("J129").Value = 2530.75
$ws.Range("K129").Value = 5462.4
$ws.Range("L129").Value = 7592.25
$ws.Range("M129").Value = -462.3999999999996
$ws.Range("N129").Value = -17592.25
$ws.Range("H132").Value = 4922.1304
$ws.Range("I132").Value = 5014.0454
$ws.Range("J132").Value = 2900
$ws.Range("K132").Value = 15042.1362
$ws.Range("L132").Value = 8700
$ws.Range("M132").Value = -12512.1362
$ws.Range("N132").Value = -13760
$ws.Range("H141").Value = 19756.125
$ws.Range("I141").Value = 29747.8
$ws.Range("K141").Value = 89243.39999999999
$ws.Range("M141").Value = -84063.39999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 18000
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H61").Value = 18520160
$ws.Range("I61").Value = 18520160
$ws.Range("K61").Value = 18520160
$ws.Range("M61").Value = -18519948
$ws.Range("H63").Value = 6853.778
$ws.Range("J63").Value = 9166.666999999999
$ws.Range("L63").Value = 9166.666999999999
$ws.Range("N63").Value = -10538.667
$ws.Range("H66").Value = 6853.778
$ws.Range("J66").Value = 9166.666999999999
$ws.Range("L66").Value = 45833.335
$ws.Range("N66").Value = -52697.335
$ws.Range("H122").Value = 8134667
$ws.Range("I122").Value = 1980.1428
$ws.Range("K122").Value = 5940.428400000001
$ws.Range("M122").Value = -3490.428400000001
$ws.Range("H132").Value = 83356710
$ws.Range("I132").Value = 14190
$ws.Range("J132").Value = 333384300
$ws.Range("K132").Value = 42570
$ws.Range("L132").Value = 1000152900
$ws.Range("M132").Value = -40040
$ws.Range("N132").Value = -1000157960
$ws.Range("H136").Value = 18520160
$ws.Range("I136").Value = 18520160
$ws.Range("K136").Value = 55560480
$ws.Range("M136").Value = -55557930

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H87").Value = 111750
$ws.Range("J87").Value = 110000
$ws.Range("L87").Value = 110000
$ws.Range("N87").Value = -112496
$ws.Range("H90").Value = 111750
$ws.Range("J90").Value = 110000
$ws.Range("L90").Value = 330000
$ws.Range("N90").Value = -342480
$ws.Range("H105").Value = 16272.571
$ws.Range("I105").Value = 27052
$ws.Range("J105").Value = 1900
$ws.Range("K105").Value = 27052
$ws.Range("L105").Value = 1900
$ws.Range("M105").Value = -25305
$ws.Range("N105").Value = -5394
$ws.Range("H134").Value = 2632.5789
$ws.Range("I134").Value = 2623.2778
$ws.Range("J134").Value = 2800
$ws.Range("K134").Value = 7869.8334
$ws.Range("L134").Value = 8400
$ws.Range("M134").Value = -5334.8334
$ws.Range("N134").Value = -13470

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 33786100
$ws.Range("I31").Value = 3526.8572
$ws.Range("J31").Value = 41668700
$ws.Range("K31").Value = 3526.8572
$ws.Range("L31").Value = 41668700
$ws.Range("M31").Value = -3231.8572
$ws.Range("N31").Value = -41669290
$ws.Range("H34").Value = 33786100
$ws.Range("I34").Value = 3526.8572
$ws.Range("J34").Value = 41668700
$ws.Range("K34").Value = 3526.8572
$ws.Range("L34").Value = 41668700
$ws.Range("M34").Value = -3324.8572
$ws.Range("N34").Value = -41669104
$ws.Range("H58").Value = 3701.6428
$ws.Range("I58").Value = 4493.6665
$ws.Range("J58").Value = 2276
$ws.Range("K58").Value = 4493.6665
$ws.Range("L58").Value = 2276
$ws.Range("M58").Value = -4290.6665
$ws.Range("N58").Value = -2682
$ws.Range("H62").Value = 1466.6666
$ws.Range("I62").Value = 1700
$ws.Range("J62").Value = 1000
$ws.Range("K62").Value = 1700
$ws.Range("L62").Value = 1000
$ws.Range("M62").Value = -1076
$ws.Range("N62").Value = -2248
$ws.Range("H65").Value = 1466.6666
$ws.Range("I65").Value = 1700
$ws.Range("J65").Value = 1000
$ws.Range("K65").Value = 8500
$ws.Range("L65").Value = 5000
$ws.Range("M65").Value = -5380
$ws.Range("N65").Value = -11240
$ws.Range("H122").Value = 2394826.2
$ws.Range("I122").Value = 2509.7778
$ws.Range("K122").Value = 7529.3334
$ws.Range("M122").Value = -5079.3334
$ws.Range("H136").Value = 3701.6428
$ws.Range("I136").Value = 4493.6665
$ws.Range("J136").Value = 2276
$ws.Range("K136").Value = 13480.9995
$ws.Range("L136").Value = 6828
$ws.Range("M136").Value = -10930.9995
$ws.Range("N136").Value = -11928

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H70").Value = 2402
$ws.Range("J70").Value = 3500
$ws.Range("L70").Value = 10500
$ws.Range("N70").Value = -11130
$ws.Range("H73").Value = 2402
$ws.Range("J73").Value = 3500
$ws.Range("L73").Value = 10500
$ws.Range("N73").Value = -12684
$ws.Range("H129").Value = 2265.7036
$ws.Range("I129").Value = 808.3333
$ws.Range("K129").Value = 2424.9999
$ws.Range("M129").Value = 2575.0001
$ws.Range("H131").Value = 2785.318
$ws.Range("I131").Value = 962
$ws.Range("K131").Value = 2886
$ws.Range("M131").Value = 2154
$ws.Range("H137").Value = 2769.9
$ws.Range("I137").Value = 1450
$ws.Range("J137").Value = 3099.875
$ws.Range("K137").Value = 4350
$ws.Range("L137").Value = 9299.625
$ws.Range("M137").Value = 750
$ws.Range("N137").Value = -19499.625
$ws.Range("H140").Value = 1165.6666
$ws.Range("I140").Value = 932.25
$ws.Range("K140").Value = 2796.75
$ws.Range("M140").Value = 2383.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10666
$ws.Range("I80").Value = 12935.571
$ws.Range("J80").Value = 6694.25
$ws.Range("K80").Value = 12935.571
$ws.Range("L80").Value = 6694.25
$ws.Range("M80").Value = -11937.571
$ws.Range("N80").Value = -8690.25
$ws.Range("H83").Value = 10666
$ws.Range("I83").Value = 12935.571
$ws.Range("J83").Value = 6694.25
$ws.Range("K83").Value = 64677.855
$ws.Range("L83").Value = 33471.25
$ws.Range("M83").Value = -59685.855
$ws.Range("N83").Value = -43455.25
$ws.Range("H132").Value = 11638.634
$ws.Range("I132").Value = 7802.6816
$ws.Range("J132").Value = 22187.5
$ws.Range("K132").Value = 23408.0448
$ws.Range("L132").Value = 66562.5
$ws.Range("M132").Value = -20878.0448
$ws.Range("N132").Value = -71622.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7329.6
$ws.Range("I40").Value = 7299.36
$ws.Range("J40").Value = 7480.8
$ws.Range("K40").Value = 7299.36
$ws.Range("L40").Value = 7480.8
$ws.Range("M40").Value = -7163.36
$ws.Range("N40").Value = -7752.8
$ws.Range("H46").Value = 1792.5714
$ws.Range("I46").Value = 990
$ws.Range("K46").Value = 990
$ws.Range("M46").Value = -802
$ws.Range("H68").Value = 6233.1665
$ws.Range("J68").Value = 7350
$ws.Range("L68").Value = 7350
$ws.Range("N68").Value = -8848
$ws.Range("H71").Value = 6233.1665
$ws.Range("J71").Value = 7350
$ws.Range("L71").Value = 36750
$ws.Range("N71").Value = -44238
$ws.Range("H82").Value = 1248.5
$ws.Range("I82").Value = 1248.5
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 1248.5
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -887.5
$ws.Range("N82").ClearContents()
$ws.Range("H85").Value = 1248.5
$ws.Range("I85").Value = 1248.5
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 1248.5
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -0.5
$ws.Range("N85").ClearContents()
$ws.Range("H93").Value = 1265824.4
$ws.Range("I93").Value = 2197
$ws.Range("K93").Value = 2197
$ws.Range("M93").Value = -949
$ws.Range("H100").Value = 2871.2
$ws.Range("I100").Value = 2495.8823
$ws.Range("K100").Value = 2495.8823
$ws.Range("M100").Value = -1954.8823
$ws.Range("H122").Value = 10420817
$ws.Range("I122").Value = 4059.9
$ws.Range("K122").Value = 12179.7
$ws.Range("M122").Value = -9729.700000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 21750
$ws.Range("J62").Value = 21750
$ws.Range("L62").Value = 21750
$ws.Range("N62").Value = -22998
$ws.Range("H65").Value = 21750
$ws.Range("J65").Value = 21750
$ws.Range("L65").Value = 108750
$ws.Range("N65").Value = -114990
$ws.Range("H81").Value = 3690.8
$ws.Range("J81").Value = 6059.8
$ws.Range("L81").Value = 12119.6
$ws.Range("N81").Value = -14241.6
$ws.Range("H84").Value = 3690.8
$ws.Range("J84").Value = 6059.8
$ws.Range("L84").Value = 60598
$ws.Range("N84").Value = -71206
$ws.Range("H122").Value = 8743132
$ws.Range("I122").Value = 54026.85
$ws.Range("K122").Value = 162080.55
$ws.Range("M122").Value = -159630.55
$ws.Range("H135").Value = 11198741
$ws.Range("J135").Value = 11198741
$ws.Range("L135").Value = 11198741
$ws.Range("N135").Value = -11208881
